$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Semestre ideal" from EA-7 to EA-8
$ws.Range("B9").Value = "EA-8"
$ws.Range("C9").Value = "EA-8"

# Update first requirement (row 24): LOB1207 Poluicao Ambiental I -> LOB1212 Quimica Analitica Ambiental II (weak requirement)
$req1 = "LOB1212 -  Química Analítica Ambiental II  (Requisito fraco)`n"
$ws.Range("B24").Value = $req1
$ws.Range("C24").Value = $req1

# Remove the second requirement row entirely (row 25: LOB1208 Quimica Analitica Ambiental I), shifting rows up
$ws.Rows("25").Delete()

# The former row 26 (LOB1255 Hidrologia Aplicada) is now row 25; update its text
$req3 = "LOB1258 -  Hidráulica Aplicada  (Requisito fraco)`n"
$ws.Range("B25").Value = $req3
$ws.Range("C25").Value = $req3
